$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.6630574638774661
$ws.Range("J2").Value = 0.663057463877466
$ws.Range("M2").Value = 0.092903
$ws.Range("N2").Value = 0.278709
$ws.Range("O2").Value = 0.03600043090620505
$ws.Range("P2").Value = 0.03600043090620505
$ws.Range("Q2").Value = 0.01514455157733333
$ws.Range("R2").Value = 0.136300964196
$ws.Range("S2").Value = 0.02387035441516427
$ws.Range("T2").Value = 0.02387035441516426
$ws.Range("I3").Value = 0.6630574638774661
$ws.Range("J3").Value = 0.663057463877466
$ws.Range("O3").Value = 0.1132051051535142
$ws.Range("P3").Value = 0.1132051051535142
$ws.Range("S3").Value = 0.07506148992107099
$ws.Range("T3").Value = 0.07506148992107101
$ws.Range("I4").Value = 0.6630574638774661
$ws.Range("J4").Value = 0.663057463877466
$ws.Range("M4").Value = 2.195567
$ws.Range("N4").Value = 6.586701
$ws.Range("O4").Value = 0.8507944639402807
$ws.Range("P4").Value = 0.8507944639402808
$ws.Range("Q4").Value = 0.3579096226493333
$ws.Range("R4").Value = 3.221186603844
$ws.Range("S4").Value = 0.5641256195412309
$ws.Range("T4").Value = 0.5641256195412307
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.08283833333333333
$ws.Range("H5").Value = 0.248515
$ws.Range("I5").Value = 0.3369425361225339
$ws.Range("J5").Value = 0.3369425361225339
$ws.Range("M5").Value = 0.092903
$ws.Range("N5").Value = 0.278709
$ws.Range("O5").Value = 0.03600043090620505
$ws.Range("P5").Value = 0.03600043090620505
$ws.Range("Q5").Value = 0.007695929681666666
$ws.Range("R5").Value = 0.06926336713499999
$ws.Range("S5").Value = 0.01213007649104078
$ws.Range("T5").Value = 0.01213007649104078
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.08283833333333333
$ws.Range("H6").Value = 0.248515
$ws.Range("I6").Value = 0.3369425361225339
$ws.Range("J6").Value = 0.3369425361225339
$ws.Range("O6").Value = 0.1132051051535142
$ws.Range("P6").Value = 0.1132051051535142
$ws.Range("Q6").Value = 0.02420022502333333
$ws.Range("R6").Value = 0.21780202521
$ws.Range("S6").Value = 0.03814361523244321
$ws.Range("T6").Value = 0.03814361523244322
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.08283833333333333
$ws.Range("H7").Value = 0.248515
$ws.Range("I7").Value = 0.3369425361225339
$ws.Range("J7").Value = 0.3369425361225339
$ws.Range("M7").Value = 2.195567
$ws.Range("N7").Value = 6.586701
$ws.Range("O7").Value = 0.8507944639402807
$ws.Range("P7").Value = 0.8507944639402808
$ws.Range("Q7").Value = 0.1818771110016667
$ws.Range("R7").Value = 1.636893999015
$ws.Range("S7").Value = 0.2866688443990499
$ws.Range("T7").Value = 0.2866688443990499
